$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that originally sat right
#    under the H1 title (paragraph #2).
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Insert a new bold paragraph "Play Bomber Squad Free: Unique 5-Reel Slot
#    with Progressive Jackpot" right before the final (italic) paragraph.
#
#    We build it right after a plain body paragraph first (so it inherits no
#    stray paragraph style/indent/italic formatting), then cut/paste it into
#    its final position right before the last paragraph. The anchor
#    paragraph is located by its (unique) trailing sentence rather than a
#    fixed index, since the deletion above shifts paragraph numbering.
# ---------------------------------------------------------------------------
$anchorText = "Who knows? You may end up blowing the vault wide open and walking away with a bag full of cash!"

$find1 = $d.Content
$find1.Find.Execute($anchorText) | Out-Null
$anchorPara = $find1.Paragraphs(1)
$anchorPara.Range.InsertParagraphAfter() | Out-Null

$find2 = $d.Content
$find2.Find.Execute($anchorText) | Out-Null
$newPara = $find2.Paragraphs(1).Next()

$insertionPoint = $newPara.Range.End - 1
$collapsed = $d.Range($insertionPoint, $insertionPoint)
$collapsed.InsertAfter("Play Bomber Squad Free: Unique 5-Reel Slot with Progressive Jackpot")

$find3 = $d.Content
$find3.Find.Execute($anchorText) | Out-Null
$newPara = $find3.Paragraphs(1).Next()
$textRange = $d.Range($insertionPoint, $newPara.Range.End - 1)
$textRange.Font.Bold = $true

$find4 = $d.Content
$find4.Find.Execute($anchorText) | Out-Null
$newPara = $find4.Paragraphs(1).Next()
$newPara.Range.Cut() | Out-Null

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$pasteRange = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$pasteRange.Paste() | Out-Null

# ---------------------------------------------------------------------------
# 3. Update the text of the final (italic) paragraph from the old image
#    prompt to the meta-description copy.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a feature image that captures the excitement and adventure of Bomber Squad by SimplePlay. The image should be in a cartoon style and should feature a happy Maya warrior wearing glasses. The Maya warrior should be holding a dynamite stick and standing in front of a bank vault door while gangsters shoot their machine guns in the background. Use bright colors and bold outlines to make the image stand out and convey the high-energy and thrilling experience of playing this slot game. The image should make players curious and excited to try out Bomber Squad.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read a review of Bomber Squad, a five-reel online slot with 1024 ways to win, a thrilling storyline, and a unique bank robbery theme. Play for free now!",
    2
) | Out-Null
